$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename sheets to reflect clearer keyword naming.
# ---------------------------------------------------------------
$wb.Worksheets.Item("Enemies").Name = "Enemy Options"
$wb.Worksheets.Item("Status Groups").Name = "Priviledged or Persecuted"

# ---------------------------------------------------------------
# 2. Military sheet: add a new "Basic Guard only" option above the
#    existing list (new row 2), pushing the rest down, matching the
#    layout/number style already used in the table.
# ---------------------------------------------------------------
$wsMilitary = $wb.Worksheets.Item("Military")
$wsMilitary.Activate()

$wsMilitary.Rows.Item(2).Insert()
$wsMilitary.Range("A2").Value = 0
$wsMilitary.Range("B2").Value = "Basic Guard only"

# Match formatting of the surrounding cells (number column / text column).
$wsMilitary.Range("A3").Copy()
$wsMilitary.Range("A2").PasteSpecial(-4122)
$wsMilitary.Range("B1").Copy()
$wsMilitary.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 3. Restore / update each sheet's selected cell where it moved.
# ---------------------------------------------------------------
$wsDiplomatic = $wb.Worksheets.Item("Diplomatic Events")
$wsDiplomatic.Activate()
$wsDiplomatic.Range("D12").Select()

$wsMagic = $wb.Worksheets.Item("Magic Area or School")
$wsMagic.Activate()
$wsMagic.Range("C20").Select()

# ---------------------------------------------------------------
# 4. Leave "Military" as the active sheet/tab, with its new
#    selection, matching the final saved workbook state.
# ---------------------------------------------------------------
$wsMilitary.Activate()
$wsMilitary.Range("B3").Select()
